$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Text='26.908.05'},
    @{Cell='E2'; Text='  +1.83%  '},
    @{Cell='D3'; Text='1.727.54'},
    @{Cell='E3'; Text='  +0.19%  '},
    @{Cell='D4'; Text='0.9967'},
    @{Cell='E4'; Text='  -0.26%  '},
    @{Cell='D5'; Text='242.26'},
    @{Cell='E5'; Text='  -0.22%  '},
    @{Cell='D6'; Text='0.9973'},
    @{Cell='E6'; Text='  -0.23%  '},
    @{Cell='D7'; Text='0.4886'},
    @{Cell='E7'; Text='  -0.61%  '},
    @{Cell='D8'; Text='0.2594'},
    @{Cell='E8'; Text='  -0.95%  '},
    @{Cell='D9'; Text='0.06213'},
    @{Cell='E9'; Text='  +0.13%  '},
    @{Cell='D10'; Text='1.731.81'},
    @{Cell='E10'; Text='  +1.17%  '},
    @{Cell='D11'; Text='16.03'},
    @{Cell='E11'; Text='  +3.38%  '},
    @{Cell='D12'; Text='0.06896'},
    @{Cell='E12'; Text='  -1.80%  '},
    @{Cell='D13'; Text='0.6091'},
    @{Cell='E13'; Text='  +1.42%  '},
    @{Cell='D14'; Text='4.489'},
    @{Cell='E14'; Text='  -1.92%  '},
    @{Cell='D15'; Text='77.23'},
    @{Cell='E15'; Text='  -0.13%  '},
    @{Cell='D16'; Text='0.9979'},
    @{Cell='E16'; Text='  -0.20%  '},
    @{Cell='D17'; Text='26.651.31'},
    @{Cell='E17'; Text='  +0.86%  '},
    @{Cell='D18'; Text='0.9965'},
    @{Cell='E18'; Text='  -0.29%  '},
    @{Cell='D19'; Text='0.000007187'},
    @{Cell='E19'; Text='  -0.03%  '},
    @{Cell='D20'; Text='11.46'},
    @{Cell='E20'; Text='  +0.84%  '},
    @{Cell='D21'; Text='1.952.65'},
    @{Cell='E21'; Text='  +0.54%  '},
    @{Cell='D22'; Text='4.423'},
    @{Cell='E22'; Text='  -1.49%  '},
    @{Cell='D23'; Text='8.587'},
    @{Cell='E23'; Text='  -0.08%  '},
    @{Cell='D24'; Text='5.096'},
    @{Cell='E24'; Text='  -1.50%  '},
    @{Cell='D25'; Text='138.63'},
    @{Cell='E25'; Text='  +0.61%  '},
    @{Cell='D26'; Text='15.28'},
    @{Cell='E26'; Text='  +0.21%  '},
    @{Cell='D27'; Text='1.775'},
    @{Cell='E27'; Text='  +3.37%  '},
    @{Cell='D28'; Text='106.28'},
    @{Cell='E28'; Text='  -0.81%  '},
    @{Cell='D29'; Text='1.380'},
    @{Cell='E29'; Text='  -1.27%  '},
    @{Cell='D30'; Text='3.955'},
    @{Cell='E30'; Text='  +0.13%  '},
    @{Cell='D31'; Text='0.08000'},
    @{Cell='E31'; Text='  +0.28%  '},
    @{Cell='D32'; Text='3.693'},
    @{Cell='E32'; Text='  +0.56%  '},
    @{Cell='D33'; Text='0.04521'},
    @{Cell='E33'; Text='  -0.45%  '},
    @{Cell='B34'; Text='Frax'},
    @{Cell='C34'; Text='https://coinranking.com/coin/KfWtaeV1W+frax-frax'},
    @{Cell='D34'; Text='0.9967'},
    @{Cell='E34'; Text='  -0.26%  '},
    @{Cell='B35'; Text='HuobiToken'},
    @{Cell='C35'; Text='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'},
    @{Cell='D35'; Text='2.596'},
    @{Cell='E35'; Text='  -0.24%  '},
    @{Cell='B36'; Text='ARBITRUM'},
    @{Cell='C36'; Text='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'},
    @{Cell='D36'; Text='1.010'},
    @{Cell='E36'; Text='  +1.39%  '},
    @{Cell='B37'; Text='ImmutableX'},
    @{Cell='C37'; Text='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'},
    @{Cell='D37'; Text='0.6260'},
    @{Cell='E37'; Text='  -0.06%  '},
    @{Cell='B38'; Text='TrustWalletToken'},
    @{Cell='C38'; Text='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'},
    @{Cell='D38'; Text='0.9371'},
    @{Cell='E38'; Text='  +1.85%  '},
    @{Cell='B39'; Text='RenderToken'},
    @{Cell='C39'; Text='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'},
    @{Cell='D39'; Text='2.052'},
    @{Cell='E39'; Text='  +4.49%  '},
    @{Cell='B40'; Text='MXToken'},
    @{Cell='C40'; Text='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Cell='D40'; Text='2.456'},
    @{Cell='E40'; Text='  +2.62%  '},
    @{Cell='B41'; Text='PaxDollar'},
    @{Cell='C41'; Text='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'},
    @{Cell='D41'; Text='0.9969'},
    @{Cell='E41'; Text='  -0.27%  '},
    @{Cell='B42'; Text='FraxShare'},
    @{Cell='C42'; Text='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'},
    @{Cell='D42'; Text='5.714'},
    @{Cell='E42'; Text='  +6.85%  '},
    @{Cell='B43'; Text='VeChain'},
    @{Cell='C43'; Text='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'},
    @{Cell='D43'; Text='0.01505'},
    @{Cell='E43'; Text='  +1.09%  '},
    @{Cell='B44'; Text='Quant'},
    @{Cell='C44'; Text='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'},
    @{Cell='D44'; Text='99.75'},
    @{Cell='E44'; Text='  -0.28%  '},
    @{Cell='B45'; Text='TheSandbox'},
    @{Cell='C45'; Text='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'},
    @{Cell='D45'; Text='0.3857'},
    @{Cell='E45'; Text='  +0.19%  '},
    @{Cell='B46'; Text='Aptos'},
    @{Cell='C46'; Text='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'},
    @{Cell='D46'; Text='6.875'},
    @{Cell='E46'; Text='  +2.20%  '},
    @{Cell='B47'; Text='Algorand'},
    @{Cell='C47'; Text='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'},
    @{Cell='D47'; Text='0.1160'},
    @{Cell='E47'; Text='  -0.51%  '},
    @{Cell='B48'; Text='Cronos'},
    @{Cell='C48'; Text='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'},
    @{Cell='D48'; Text='0.05401'},
    @{Cell='E48'; Text='  +0.69%  '},
    @{Cell='B49'; Text='EnergySwap'},
    @{Cell='C49'; Text='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
    @{Cell='D49'; Text='7.901'},
    @{Cell='E49'; Text='  +2.50%  '},
    @{Cell='B50'; Text='Elrond'},
    @{Cell='C50'; Text='https://coinranking.com/coin/omwkOTglq+elrond-egld'},
    @{Cell='D50'; Text='30.28'},
    @{Cell='E50'; Text='  +0.44%  '},
    @{Cell='B51'; Text='Aave'},
    @{Cell='C51'; Text='https://coinranking.com/coin/ixgUfzmLR+aave-aave'},
    @{Cell='D51'; Text='51.74'},
    @{Cell='E51'; Text='  +1.57%  '}
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $t = $u.Text
    if ($t -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $c.NumberFormat = "@"
        $c.Value = $t
        $c.Style = "Normal"
    } else {
        $c.Value = $t
    }
}

